# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.740.38'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.821.28'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.28'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.581'
$ws.Range('E6').Value = '  +4.46%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.98'
$ws.Range('E8').Value = '  +7.24%  '
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0700'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0952'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '2.085.27'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.40'
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('D14').Value = '1.827.57'
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.647'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').Value = '34.703.58'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.35'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.38'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '0.0₃0802'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '246.59'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.63'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '173.68'
$ws.Range('E24').Value = '  +6.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.11'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.52'
$ws.Range('E26').Value = '  +3.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.90'
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.119'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.99'
$ws.Range('E30').Value = '  +2.43%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0531'
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.25'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.404.61'
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.57'
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.679'
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0192'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.88'
$ws.Range('E40').Value = '  +4.84%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '83.78'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.79'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0516'
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.05'
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('D48').Value = '1.984.81'
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.34'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('E51').Value = '  +0.04%  '
